$wb = $excel.ActiveWorkbook

# Sheet "max-arrecad": reorder categoria_mencao labels for rows 2-6
$wsMax = $wb.Worksheets.Item("max-arrecad")
$wsMax.Range("A2").Value = "humor"
$wsMax.Range("A3").Value = "religiosidade"
$wsMax.Range("A4").Value = "terror"
$wsMax.Range("A5").Value = "hqmix"
$wsMax.Range("A6").Value = "angelo_agostini"

# Sheet "tx-sucesso": swap categoria_mencao labels for rows 15-16
$wsTx = $wb.Worksheets.Item("tx-sucesso")
$wsTx.Range("A15").Value = "politica"
$wsTx.Range("A16").Value = "erotismo"
